$wb = $excel.ActiveWorkbook

# --- Add the new "indications_interventions" sheet at the end of the tab strip ---
# The target sheetId is 6 (one higher than would be assigned if this were the
# first new sheet added in this session), so we create a throw-away sheet
# first to consume sheetId 5 / rId5, then delete it once the real sheet exists.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tempSheet = $wb.Worksheets.Add($null, $lastSheet)
$tempSheet.Name = "zzz_temp_placeholder"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet2)
$newSheet.Name = "indications_interventions"

$wb.Worksheets.Item("zzz_temp_placeholder").Delete()

# Re-fetch the new sheet by name: the old object reference can go stale once
# another sheet has been deleted, so look it up fresh before using it again.
$newSheet = $wb.Worksheets.Item("indications_interventions")

# --- Populate the new sheet, column by column then row by row, matching the
#     original authoring order so shared-string indices line up ---
$newSheet.Range("A1").Value = "type"
$newSheet.Range("A2").Value = "IND"
$newSheet.Range("A3").Value = "INT"
$newSheet.Range("A4").Value = "IND"
$newSheet.Range("A5").Value = "INT"

$newSheet.Range("B1").Value = "description"
$newSheet.Range("C1").Value = "codes"

$newSheet.Range("B2").Value = "An indication"
$newSheet.Range("C2").Value = "'SNOMED:12345=Indication1"

$newSheet.Range("B3").Value = "An intervention"
$newSheet.Range("C3").Value = "ICD-10: X = Y, SNOMED: A=B"

$newSheet.Range("B4").Value = "An indication"
$newSheet.Range("C4").Value = "'SNOMED:345678=Indication2"

$newSheet.Range("B5").Value = "An intervention"
$newSheet.Range("C5").Value = "ICD-10: DD=CC, SNOMED: A=B"

# Header row is bold (reuses the same style as the other sheets' header rows)
$newSheet.Range("A1:C1").Font.Bold = $true

# Column widths (chosen so the saved column width matches the target as
# closely as the engine's internal rounding allows)
$newSheet.Columns.Item(2).ColumnWidth = 18.166666666666668
$newSheet.Columns.Item(3).ColumnWidth = 49.5

# View state for the new sheet: zoomed to 160%, C14 selected, and it becomes
# the active/selected tab
$newSheet.Range("C14").Select()
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 160

# --- Update the previously-active "study" sheet's view state: it is no
#     longer the selected tab, and its selection becomes the full used range ---
$studySheet = $wb.Worksheets.Item("study")
$studySheet.Range("A1:D2").Select()

# Re-activate the new sheet so it ends up as the workbook's active tab
$wb.Worksheets.Item("indications_interventions").Activate()
